$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.989.13"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.827.65"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "0.9958"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "241.14"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "0.6164"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").Value = "0.9983"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.07429"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "0.2914"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "22.99"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "0.07635"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "1.825.83"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "4.981"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "0.6711"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "82.62"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "0.000009147"
$ws.Range("D17").Value = "5.883"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "28.961.07"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "2.074.60"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "239.85"
$ws.Range("E20").Value = "  +6.44%  "
$ws.Range("D21").Value = "12.63"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "0.9981"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "7.177"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "0.9970"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "158.25"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").Value = "0.1404"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "8.474"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "17.82"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "1.491"
$ws.Range("D30").Value = "0.05589"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("D31").Value = "4.129"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "4.101"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "1.197"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "1.837"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "0.7390"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "1.138"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "2.650"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").Value = "2.765"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "0.01782"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "1.210.42"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("D41").Value = "6.406"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").Value = "0.8941"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "0.9972"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "101.17"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "1.974.01"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").Value = "65.21"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "0.5069"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "0.4049"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "9.120"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "0.05800"
$ws.Range("E51").Value = "  +0.31%  "
